# Working schedule update: add checkmarks to a handful of cells in the
# "Gerçekleşen Proje Zaman Çizelgesi" block (columns S/T) and move the
# active selection from N5 to M5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$check = [char]0x2713  # "✓" already present in sharedStrings (index 13)

# Use a cell that already carries the target "checkmark" cell style
# (border + centered Times New Roman) as the format source, so the
# copied style resolves to the same existing cellXfs entry instead of
# creating a new one.
$styleSource = $ws.Range("D4")

$checkCells = @("S4", "T4", "S6", "S8", "T9", "T10", "T11")
foreach ($addr in $checkCells) {
    $styleSource.Copy()
    $ws.Range($addr).PasteSpecial(-4122) # xlPasteFormats
    $ws.Range($addr).Value = $check
}

# T12 only changes style (to the same style as the checkmark cells)
# but stays empty - no value assigned.
$styleSource.Copy()
$ws.Range("T12").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Move the active selection from N5 to M5.
$ws.Range("M5").Select()
